# Updated reliability & robustness 11/16/2016
$wb = $excel.ActiveWorkbook

$wsEff = $wb.Worksheets.Item("Efficiency")
$wsRel = $wb.Worksheets.Item("Reliability")
$wsRob = $wb.Worksheets.Item("Robustness")

$newDate = [datetime]"2016-11-16"

# =======================================================================
# Efficiency sheet: the existing dated-heading style (row 3 / C3) becomes
# bold, and a new (still empty) dated-heading row is added at row 10,
# ready for the next data drop.
# =======================================================================
$wsEff.Range("C3").Font.Bold = $true

$wsEff.Range("C10").NumberFormat = $wsEff.Range("C3").NumberFormat
$wsEff.Range("C10").Font.Bold = $true
$wsEff.Range("C10").Font.Size = 18
$wsEff.Rows.Item(10).RowHeight = 23.25

# =======================================================================
# Reliability sheet: a second Travis-CI build-time sample (11/16/2016) is
# added in columns G:J, mirroring the existing A:D table, and the
# AVERAGE formula for the first table moves from E7 down to D8 (directly
# under the header row, matching the new table's layout).
# =======================================================================
$avgFormula = $wsRel.Range("E7").Formula
$wsRel.Range("E7").Clear() | Out-Null

$wsRel.Range("D8").Formula = $avgFormula
$wsRel.Range("D8").NumberFormat = $wsRel.Range("B8").NumberFormat

$wsRel.Range("I3").Value = $newDate
$wsRel.Range("I3").NumberFormat = $wsRel.Range("C3").NumberFormat
$wsRel.Range("I3").Font.Size = $wsRel.Range("C3").Font.Size
$wsRel.Range("I3").Font.Bold = $wsRel.Range("C3").Font.Bold

$wsRel.Range("G7").Value2 = $wsRel.Range("A7").Value2
$wsRel.Range("H7").Value2 = $wsRel.Range("B7").Value2
$wsRel.Range("I7").Value2 = $wsRel.Range("C7").Value2
$wsRel.Range("J7").Value2 = $wsRel.Range("D7").Value2
$wsRel.Range("G7:J7").Font.Bold = $true
$wsRel.Range("H7").NumberFormat = $wsRel.Range("B7").NumberFormat

$wsRel.Range("G10").Value2 = "#176"
$wsRel.Range("H10").Value2 = 653
$wsRel.Range("G9").Value2 = "#209"
$wsRel.Range("H9").Value2 = 49
$wsRel.Range("G8").Value2 = "#225"
$wsRel.Range("H8").Value2 = 47

$wsRel.Range("J8").Formula = "=AVERAGE(H8:H69)"
$wsRel.Range("J8").NumberFormat = $wsRel.Range("D8").NumberFormat

$wsRel.Columns.Item(7).ColumnWidth = 21.5703125
$wsRel.Columns.Item(8).ColumnWidth = 16.140625
$wsRel.Columns.Item(9).ColumnWidth = 15.28515625
$wsRel.Columns.Item(10).ColumnWidth = 30.7109375

# =======================================================================
# Robustness sheet: a second bug-risk sample (11/16/2016) is added in
# column D, mirroring the existing column A/B figures.
# =======================================================================
$wsRob.Range("D3").Value = $newDate
$wsRob.Range("D3").NumberFormat = $wsRob.Range("A3").NumberFormat

$wsRob.Range("D5").Value2 = "43 / 52"

$wsRob.Range("D6").Value2 = 0.82
$wsRob.Range("D6").NumberFormat = $wsRob.Range("B6").NumberFormat

$wsRob.Columns.Item(4).ColumnWidth = 15.42578125

# =======================================================================
# Selections + active tab: Robustness becomes the active (last-viewed)
# sheet, with each sheet's selection parked at its newest cell.
# =======================================================================
$wsEff.Range("C10").Select() | Out-Null
$wsRel.Range("I3").Select() | Out-Null
$wsRob.Range("D7").Select() | Out-Null

$wsRob.Activate()
